# Weekly update: insert two new "Ciruela" (plum) price records at the top
# of the data table (rows 23-24), pushing the existing rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 23; Excel shifts rows 23:83 down to
# 25:85 and the sheet dimension grows from A1:T83 to A1:T85 automatically.
$ws.Rows("23:24").Insert()

# New row 23: Black Amber / Especial
$ws.Range("A23").Value = 9
$ws.Range("B23").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C23").Value = "Metropolitana"
$ws.Range("D23").Value = 44592
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100103
$ws.Range("H23").Value = "Frutos de hueso (carozo)"
$ws.Range("I23").Value = 100103002
$ws.Range("J23").Value = "Ciruela"
$ws.Range("K23").Value = "Black Amber"
$ws.Range("L23").Value = "Especial"
$ws.Range("M23").Value = 380
$ws.Range("N23").Value = 11000
$ws.Range("O23").Value = 12000
$ws.Range("P23").Value = 11474
$ws.Range("Q23").Value = "`$/caja 15 kilos granel"
$ws.Range("R23").Value = "Región de O'Higgins"
$ws.Range("S23").Value = 765
$ws.Range("T23").Value = 15

# New row 24: Black Amber / Primera
$ws.Range("A24").Value = 9
$ws.Range("B24").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C24").Value = "Metropolitana"
$ws.Range("D24").Value = 44592
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100103
$ws.Range("H24").Value = "Frutos de hueso (carozo)"
$ws.Range("I24").Value = 100103002
$ws.Range("J24").Value = "Ciruela"
$ws.Range("K24").Value = "Black Amber"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 450
$ws.Range("N24").Value = 9000
$ws.Range("O24").Value = 10000
$ws.Range("P24").Value = 9556
$ws.Range("Q24").Value = "`$/caja 15 kilos granel"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 637
$ws.Range("T24").Value = 15
